$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "310.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.34%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "7"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.58%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "7"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.00%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "7"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07759"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-4.97%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "7"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.344"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.37%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "7"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.212"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.38%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "7"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.886"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.22%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "7"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.881"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-12.44%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "7"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9167"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.47%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "7"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1199"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-8.80%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "7"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1914"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.06%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "7"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08872"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.04%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "7"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03394"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.73%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "7"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09693"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.23%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "7"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001364"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.68%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "7"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005875"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.28%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "7"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.58%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "7"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.75%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "7"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.035"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.52%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "7"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1262"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-4.15%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "7"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.99%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "7"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5,589.57%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "7"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04384"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.29%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "7"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.78%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "7"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004253"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-10.42%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "7"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-65.33%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "7"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "7"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "7"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "7"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "7"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "7"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "7"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "7"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "7"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "7"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "7"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "7"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02112"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.86%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "7"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04950"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.63%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "7"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007623"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.70%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "7"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009935"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.58%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "7"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1339"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.99%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "7"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002059"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.08%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "7"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009611"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.28%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "7"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006592"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.60%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "7"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.18%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "7"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003041"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.82%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "7"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "7"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.18%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "7"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.18%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "7"
